$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: B3 becomes "A,B", C3 cleared
$ws.Range("B3").Value = "A,B"
$ws.Range("C3").ClearContents()

# Row 4: B4 cleared, C4 becomes "B,C"
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "B,C"

# Row 5: B5 becomes "C", C5 cleared, D5 becomes 1
$ws.Range("B5").Value = "C"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 1

# Row 6 removed entirely
$ws.Range("A6:D6").ClearContents()

$ws.Range("D6").Select()
